# Insert a new "test_reference" column immediately before the existing
# "runner_settings" column on the TestAsset, AcceptanceTestAsset, and
# TestEdgeData sheets, shifting runner_settings/id/name/description/tags
# one column to the right.

$wb = $excel.ActiveWorkbook

$sheetsToUpdate = @("TestAsset", "AcceptanceTestAsset", "TestEdgeData")

foreach ($sheetName in $sheetsToUpdate) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the column whose header is "runner_settings" in row 1.
    $headerCell = $ws.Rows.Item(1).Find("runner_settings")
    $col = $headerCell.Column

    # Insert a new blank column before it, shifting existing columns right.
    $ws.Columns.Item($col).Insert()

    # Populate the new column's header cell.
    $ws.Cells.Item(1, $col).Value = "test_reference"
}
